$d = $word.ActiveDocument

# Locate the anchor paragraph ("Want- get student image upload status")
# and then the paragraph two below it, which is the lone "<tab>" paragraph
# that is being expanded into the new "Student attendance" section.
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.Trim()
    if ($t -eq "Want- get student image upload status") {
        $target = $d.Paragraphs($i + 2)
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate target paragraph"
}

$xml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr></w:pPr><w:r><w:t>Student- Student get attendance data</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>page=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>student&amp;action</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>get_student_attendance</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:tab/><w:t>Request: token</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Response: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Attendance_data</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramStart"/><w:r><w:t>&lt;</w:t></w:r><w:r><w:t>{ [</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>key: string]: { [key: string]: number } }</w:t></w:r><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Want-get attendance student </w:t></w:r></w:p>'

$target.Range.InsertXML($xml)

Write-Output "Inserted attendance student section"
